$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellValues = @{
    "D2" = 7975
    "E2" = 935
    "F2" = 935
    "G2" = 712
    "H2" = 576
    "I2" = 559
    "J2" = 17
    "K2" = 10434
    "L2" = 4750
    "M2" = 5684
    "N2" = 5736
    "O2" = -52
    "P2" = 522
    "Q2" = 390
    "R2" = -466
    "S2" = 83
    "T2" = 65
    "U2" = 325
    "V2" = 2932
    "W2" = 11.73
    "X2" = 7.22
    "Y2" = 10.39
    "Z2" = 5.72
    "AA2" = 83.56999999999999
    "AB2" = 1000.97
    "AC2" = 1084
    "AD2" = 20.58
    "AE2" = 10980
    "AF2" = 2.03
    "AG2" = 50
    "AH2" = 0.22
    "AI2" = 4.67
    "AJ2" = 52240070
    "D3" = 8158
    "E3" = 806
    "F3" = 806
    "G3" = -1756
    "H3" = -1333
    "I3" = -1328
    "J3" = -5
    "K3" = 9789
    "L3" = 4998
    "M3" = 4791
    "N3" = 4799
    "O3" = -8
    "P3" = 565
    "Q3" = 274
    "R3" = -613
    "S3" = 175
    "T3" = 132
    "U3" = 142
    "V3" = 3140
    "W3" = 9.880000000000001
    "X3" = -16.34
    "Y3" = -25.21
    "Z3" = -13.19
    "AA3" = 104.33
    "AB3" = 707.92
    "AC3" = -2399
    "AD3" = -7.87
    "AE3" = 8488
    "AF3" = 2.22
    "AG3" = 50
    "AH3" = 0.26
    "AI3" = -2.13
    "AJ3" = 56540070
    "D4" = 9671
    "E4" = 118
    "F4" = 118
    "G4" = 4311
    "H4" = 3111
    "I4" = 3233
    "J4" = -122
    "K4" = 33849
    "L4" = 20828
    "M4" = 13022
    "N4" = 8576
    "O4" = 4446
    "P4" = 570
    "Q4" = -35
    "R4" = -2691
    "S4" = 3701
    "T4" = 233
    "U4" = -268
    "V4" = 11908
    "W4" = 1.22
    "X4" = 32.17
    "Y4" = 48.35
    "Z4" = 14.26
    "AA4" = 159.95
    "AB4" = 1351.13
    "AC4" = 5684
    "AD4" = 2.48
    "AE4" = 15035
    "AF4" = 0.9399999999999999
    "AG4" = 50
    "AH4" = 0.35
    "AI4" = 0.88
    "AJ4" = 57040070
    "D5" = 25303
    "E5" = 2175
    "F5" = 2175
    "G5" = 1818
    "H5" = 1081
    "I5" = 624
    "J5" = 457
    "K5" = 30789
    "L5" = 18446
    "M5" = 12343
    "N5" = 7990
    "O5" = 4353
    "P5" = 611
    "Q5" = 300
    "R5" = -516
    "S5" = 222
    "T5" = 393
    "U5" = -94
    "V5" = 11466
    "W5" = 8.59
    "X5" = 4.27
    "Y5" = 7.53
    "Z5" = 3.35
    "AA5" = 149.44
    "AB5" = 1304.79
    "AC5" = 1030
    "AD5" = 15.85
    "AE5" = 13074
    "AF5" = 1.25
    "AG5" = 50
    "AH5" = 0.31
    "AI5" = 4.9
    "AJ5" = 61115070
    "D6" = 29546
    "E6" = 3571
    "F6" = 3571
    "G6" = 3371
    "H6" = 2100
    "I6" = 1435
    "K6" = 33043
    "L6" = 17978
    "M6" = 15065
    "N6" = 9917
    "P6" = 611
    "Q6" = 2746
    "R6" = -246
    "S6" = -2144
    "T6" = 847
    "U6" = 1899
    "V6" = 9832
    "W6" = 12.09
    "X6" = 7.11
    "Y6" = 16.03
    "Z6" = 6.58
    "AA6" = 119.34
    "AB6" = 1538.42
    "AC6" = 2349
    "AD6" = 22.78
    "AE6" = 16226
    "AF6" = 3.3
    "AG6" = 50
    "AH6" = 0.09
    "AI6" = 2.13
    "AJ6" = 61115070
    "D7" = 34822
    "E7" = 4878
    "G7" = 4928
    "H7" = 3583
    "I7" = 2761
    "K7" = 37384
    "L7" = 18685
    "M7" = 18698
    "N7" = 12870
    "P7" = 610
    "Q7" = 3382
    "R7" = -702
    "S7" = -970
    "T7" = 545
    "U7" = 3066
    "W7" = 14.01
    "X7" = 10.29
    "Y7" = 24.23
    "Z7" = 10.17
    "AA7" = 99.93000000000001
    "AC7" = 4517
    "AD7" = 9.789999999999999
    "AE7" = 21159
    "AF7" = 2.09
    "AG7" = 134
    "AH7" = 0.3
    "AI7" = 2.97
    "D8" = 36822
    "E8" = 5277
    "G8" = 5273
    "H8" = 3863
    "I8" = 3016
    "K8" = 40686
    "L8" = 18401
    "M8" = 22285
    "N8" = 15853
    "P8" = 610
    "Q8" = 4152
    "R8" = -771
    "S8" = -766
    "T8" = 732
    "U8" = 2881
    "W8" = 14.33
    "X8" = 10.49
    "Y8" = 21
    "Z8" = 9.9
    "AA8" = 82.56999999999999
    "AC8" = 4935
    "AD8" = 8.960000000000001
    "AE8" = 26095
    "AF8" = 1.69
    "AG8" = 211
    "AH8" = 0.48
    "AI8" = 4.28
    "D9" = 38972
    "E9" = 5684
    "G9" = 5820
    "H9" = 4243
    "I9" = 3356
    "K9" = 44370
    "L9" = 18015
    "M9" = 26354
    "N9" = 19141
    "P9" = 610
    "Q9" = 4672
    "R9" = -877
    "S9" = -831
    "T9" = 894
    "U9" = 3260
    "W9" = 14.58
    "X9" = 10.89
    "Y9" = 19.18
    "Z9" = 9.98
    "AA9" = 68.36
    "AC9" = 5491
    "AD9" = 8.050000000000001
    "AE9" = 31508
    "AF9" = 1.4
    "AG9" = 272
    "AH9" = 0.62
    "AI9" = 4.96
}

foreach ($ref in $cellValues.Keys) {
    $ws.Range($ref).Value = $cellValues[$ref]
}
